$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 1108 (shifts existing rows 1108-1211 down to 1110-1213)
$ws.Rows("1108:1109").Insert()

# ---- New row 1108 (Betarraga, Primera, week of 2023-08-28) ----
$ws.Range("A1108").Value = 3
$ws.Range("B1108").Value = "Femacal de La Calera"
$ws.Range("C1108").Value = "Coquimbo"
$ws.Range("D1108").Value = 45166
$ws.Range("E1108").Value = 5
$ws.Range("F1108").Value = 100114014
$ws.Range("G1108").Value = "Betarraga"
$ws.Range("H1108").Value = "Sin especificar"
$ws.Range("I1108").Value = "Primera"
$ws.Range("J1108").Value = 3200
$ws.Range("K1108").Value = 550
$ws.Range("L1108").Value = 600
$ws.Range("M1108").Value = 575
$ws.Range("N1108").Value = "$/paquete 4 unidades"
$ws.Range("O1108").Value = "Provincia de Quillota"
$ws.Range("P1108").Value = 144
$ws.Range("Q1108").Value = 4
$ws.Range("R1108").Value = "Hortaliza"

# ---- New row 1109 (Betarraga, Segunda, week of 2023-08-28) ----
$ws.Range("A1109").Value = 3
$ws.Range("B1109").Value = "Femacal de La Calera"
$ws.Range("C1109").Value = "Coquimbo"
$ws.Range("D1109").Value = 45166
$ws.Range("E1109").Value = 5
$ws.Range("F1109").Value = 100114014
$ws.Range("G1109").Value = "Betarraga"
$ws.Range("H1109").Value = "Sin especificar"
$ws.Range("I1109").Value = "Segunda"
$ws.Range("J1109").Value = 1300
$ws.Range("K1109").Value = 450
$ws.Range("L1109").Value = 450
$ws.Range("M1109").Value = 450
$ws.Range("N1109").Value = "$/paquete 4 unidades"
$ws.Range("O1109").Value = "Provincia de Quillota"
$ws.Range("P1109").Value = 112
$ws.Range("Q1109").Value = 4
$ws.Range("R1109").Value = "Hortaliza"
